# första 3 GAT 4l pnorm :rocket:
$wb = $excel.ActiveWorkbook

$wsGAT = $wb.Worksheets.Item("GAT")

# --- Update the 3rd table's title (row 36) on the GAT sheet: it now
# describes the 4-layer P-norm run specifically (vs. the generic
# placeholder title used previously). ---
$wsGAT.Range("C36").Value = "P-norm, nn.Parameter(torch.rand(out_dim)*3+1) , torch.clamp(self.p,1,100), PATTERN, "

# --- Fill in the results for the first three completed seeds (41, 35, 12)
# of the new 4-layer P-norm GAT run. ---
$wsGAT.Range("C38").Value = 41
$wsGAT.Range("D38").Value = 4
$wsGAT.Range("E38").Value = 110540
$wsGAT.Range("F38").Value = 85.662999999999997
$wsGAT.Range("G38").Value = 86.120500000000007
$wsGAT.Range("H38").Value = 63
$wsGAT.Range("I38").Value = 382.9778
$wsGAT.Range("J38").Value = 24661.905999999999
$wsGAT.Range("K38").Value = "P100"

$wsGAT.Range("C39").Value = 35
$wsGAT.Range("D39").Value = 4
$wsGAT.Range("E39").Value = 110540
$wsGAT.Range("F39").Value = 85.691999999999993
$wsGAT.Range("G39").Value = 86.033900000000003
$wsGAT.Range("H39").Value = 57
$wsGAT.Range("I39").Value = 386.4939
$wsGAT.Range("J39").Value = 22573.454300000001
$wsGAT.Range("K39").Value = "P100"

$wsGAT.Range("C40").Value = 12
$wsGAT.Range("D40").Value = 4
$wsGAT.Range("E40").Value = 110540
$wsGAT.Range("F40").Value = 85.64
$wsGAT.Range("G40").Value = 85.933300000000003
$wsGAT.Range("H40").Value = 59
$wsGAT.Range("I40").Value = 288.84780000000001
$wsGAT.Range("J40").Value = 23492.045300000002
$wsGAT.Range("K40").Value = "P100"

# 4th seed (95) is still training - only SEED/Layers/Params are known so
# far. It's running on a V100 this time, flagged "Ej klar" (not done).
$wsGAT.Range("C41").Value = 95
$wsGAT.Range("D41").Value = 4
$wsGAT.Range("E41").Value = 110540
$wsGAT.Range("L41").Value = "Ej klar"
$wsGAT.Range("K41").Value = "V100"

# --- Recreate the author's final UI state: they were scrolled through
# GIN and SAGE, and ended up on GAT (which becomes the active/selected
# tab on save). ---
$wsGIN = $wb.Worksheets.Item("GIN")
$wsGIN.Activate()
$wsGIN.Range("M38").Select() | Out-Null

$wsSAGE = $wb.Worksheets.Item("SAGE")
$wsSAGE.Activate()
$wsSAGE.Range("N18").Select() | Out-Null

$wsGAT.Activate()
$wsGAT.Range("J48").Select() | Out-Null
